# Trade #18 closed at 2026-02-17 20:04:11 - unknown UNKNOWN +0.000%
#
# Updates the "live trading results" workbook after the 18th trade
# (Trade # 18, strategy "MarketMaking") closed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1399.79     # Current Capital
$summary.Range("B4").Value = -0.22       # Total P&L $
$summary.Range("B5").Value = -0.24       # Total P&L %
$summary.Range("B6").Value = 18          # Total Trades
$summary.Range("B7").Value = 9           # Winning Trades
$summary.Range("B9").Value = 50          # Win Rate %

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 99.79000000000001   # Capital
$status.Range("D5").Value = 18                  # Trades
$status.Range("E5").Value = -0.22               # P&L $
$status.Range("F5").Value = -0.21               # P&L %
$status.Range("G5").Value = 50                  # Win Rate %

# ---------------------------------------------------------------------
# 3) Append the new trade row (row 19) to "All Trades" and
#    "MarketMaking" sheets
# ---------------------------------------------------------------------
$tradeRow = @{
    A = 18
    B = "2026-02-17"
    C = "20:04:05"
    D = "MarketMaking"
    E = "UP"
    F = 0.95
    G = 0.96
    H = "CLOSED"
    I = 1.0526
    J = 0.01
    K = 99.79000000000001
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.14
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A19").Value = $tradeRow.A

    # B19 holds a "yyyy-mm-dd" looking string that must stay plain text
    # (matching every other cell in the Date column) instead of being
    # auto-converted to a date serial by the usual Value-assignment
    # heuristic. Force text formatting for the assignment, then strip the
    # formatting back off so the cell ends up styled exactly like its
    # neighbours (no explicit style index).
    $ws.Range("B19").NumberFormat = "@"
    $ws.Range("B19").Value = $tradeRow.B
    $ws.Range("B19").ClearFormats()

    $ws.Range("C19").Value = $tradeRow.C
    $ws.Range("D19").Value = $tradeRow.D
    $ws.Range("E19").Value = $tradeRow.E
    $ws.Range("F19").Value = $tradeRow.F
    $ws.Range("G19").Value = $tradeRow.G
    $ws.Range("H19").Value = $tradeRow.H
    $ws.Range("I19").Value = $tradeRow.I
    $ws.Range("J19").Value = $tradeRow.J
    $ws.Range("K19").Value = $tradeRow.K
    $ws.Range("L19").Value = $tradeRow.L
    $ws.Range("M19").Value = $tradeRow.M
    $ws.Range("N19").Value = $tradeRow.N
    $ws.Range("O19").Value = $tradeRow.O
    $ws.Range("P19").Value = $tradeRow.P
    $ws.Range("Q19").Value = $tradeRow.Q
}
